$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(4, 6).Value = 27
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(12, 6).Value = 5504
$ws.Cells.Item(13, 6).Value = 72
$ws.Cells.Item(14, 6).Value = 6285
$ws.Cells.Item(16, 6).Value = 415
$ws.Cells.Item(17, 6).Value = 410
$ws.Cells.Item(18, 6).Value = 31
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(34, 6).Value = 74
$ws.Cells.Item(36, 6).Value = 312
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(39, 6).Value = 5243
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(45, 6).Value = 1078
$ws.Cells.Item(46, 6).Value = 989
$ws.Cells.Item(47, 6).Value = 1376
$ws.Cells.Item(49, 6).Value = 1096

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 38
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(19, 6).Value = 0

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(6, 6).Value = 16
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(9, 6).Value = 148
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(11, 6).Value = 7153
$ws.Cells.Item(12, 6).Value = 185
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(16, 6).Value = 5504
$ws.Cells.Item(17, 6).Value = 72
$ws.Cells.Item(18, 6).Value = 6285
$ws.Cells.Item(19, 6).Value = 6285
$ws.Cells.Item(24, 6).Value = 276
$ws.Cells.Item(25, 6).Value = 209
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(30, 6).Value = 47
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(39, 6).Value = 1397
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(42, 6).Value = 670
$ws.Cells.Item(43, 6).Value = 115
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(48, 6).Value = 0
$ws.Cells.Item(49, 6).Value = 0
